# Auto-generated edit script: updates cached market-price / profit
# values across the Leve profit tables on several sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1677.48
$ws.Range("I15").Value = 1677.48
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5032.440000000001
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4863.440000000001
$ws.Range("H69").Value = 3249.8
$ws.Range("I69").Value = 2910
$ws.Range("J69").Value = 3509.647
$ws.Range("K69").Value = 8730
$ws.Range("L69").Value = 10528.941
$ws.Range("M69").Value = -7856
$ws.Range("N69").Value = -12276.941
$ws.Range("H72").Value = 3249.8
$ws.Range("I72").Value = 2910
$ws.Range("J72").Value = 3509.647
$ws.Range("K72").Value = 26190
$ws.Range("L72").Value = 31586.823
$ws.Range("M72").Value = -21822
$ws.Range("N72").Value = -40322.823
$ws.Range("H80").Value = 490.47058
$ws.Range("I80").Value = 593.8889
$ws.Range("J80").Value = 374.125
$ws.Range("K80").Value = 1781.6667
$ws.Range("L80").Value = 1122.375
$ws.Range("M80").Value = -783.6667000000002
$ws.Range("N80").Value = -3118.375
$ws.Range("H83").Value = 490.47058
$ws.Range("I83").Value = 593.8889
$ws.Range("J83").Value = 374.125
$ws.Range("K83").Value = 5345.0001
$ws.Range("L83").Value = 3367.125
$ws.Range("M83").Value = -353.0001000000002
$ws.Range("N83").Value = -13351.125
$ws.Range("H100").Value = 15923746
$ws.Range("I100").Value = 25643296
$ws.Range("J100").Value = 129475.25
$ws.Range("K100").Value = 25643296
$ws.Range("L100").Value = 129475.25
$ws.Range("M100").Value = -25642755
$ws.Range("N100").Value = -130557.25
$ws.Range("H141").Value = 3669.8286
$ws.Range("I141").Value = 1607.625
$ws.Range("J141").Value = 25666.666
$ws.Range("K141").Value = 4822.875
$ws.Range("L141").Value = 76999.99800000001
$ws.Range("M141").Value = 357.125
$ws.Range("N141").Value = -87359.99800000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1337.8441
$ws.Range("I61").Value = 1164.4286
$ws.Range("J61").Value = 2118.2144
$ws.Range("K61").Value = 1164.4286
$ws.Range("L61").Value = 2118.2144
$ws.Range("M61").Value = -952.4286
$ws.Range("N61").Value = -2542.2144
$ws.Range("H74").Value = 3377.9443
$ws.Range("I74").Value = 3664.8838
$ws.Range("J74").Value = 2256.2727
$ws.Range("K74").Value = 3664.8838
$ws.Range("L74").Value = 2256.2727
$ws.Range("M74").Value = -2790.8838
$ws.Range("N74").Value = -4004.2727
$ws.Range("H77").Value = 3377.9443
$ws.Range("I77").Value = 3664.8838
$ws.Range("J77").Value = 2256.2727
$ws.Range("K77").Value = 18324.419
$ws.Range("L77").Value = 11281.3635
$ws.Range("M77").Value = -13956.419
$ws.Range("N77").Value = -20017.3635
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 1337.8441
$ws.Range("I136").Value = 1164.4286
$ws.Range("J136").Value = 2118.2144
$ws.Range("K136").Value = 3493.2858
$ws.Range("L136").Value = 6354.6432
$ws.Range("M136").Value = -943.2857999999997
$ws.Range("N136").Value = -11454.6432

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 533.3333
$ws.Range("I22").Value = 533.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 533.3333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -360.3333
$ws.Range("N22").ClearContents()
$ws.Range("H26").Value = 18000
$ws.Range("I26").Value = 18000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 18000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -17708
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H96").Value = 1950
$ws.Range("I96").Value = 1950
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1950
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 796
$ws.Range("N96").ClearContents()
$ws.Range("H105").Value = 1424
$ws.Range("I105").Value = 1126
$ws.Range("J105").Value = 2020
$ws.Range("K105").Value = 1126
$ws.Range("L105").Value = 2020
$ws.Range("M105").Value = 621
$ws.Range("N105").Value = -5514

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 404.92307
$ws.Range("I22").Value = 172
$ws.Range("J22").Value = 550.5
$ws.Range("K22").Value = 172
$ws.Range("L22").Value = 550.5
$ws.Range("M22").Value = 178
$ws.Range("N22").Value = -1250.5
$ws.Range("H105").Value = 1715.0416
$ws.Range("I105").Value = 2232.1428
$ws.Range("J105").Value = 991.1
$ws.Range("K105").Value = 2232.1428
$ws.Range("L105").Value = 991.1
$ws.Range("M105").Value = -485.1428000000001
$ws.Range("N105").Value = -4485.1
$ws.Range("H134").Value = 1481.1061
$ws.Range("I134").Value = 1413.4807
$ws.Range("J134").Value = 1732.2858
$ws.Range("K134").Value = 4240.4421
$ws.Range("L134").Value = 5196.857400000001
$ws.Range("M134").Value = -1705.4421
$ws.Range("N134").Value = -10266.8574

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 439323.9
$ws.Range("I5").Value = 457
$ws.Range("J5").Value = 798396.8
$ws.Range("K5").Value = 1371
$ws.Range("L5").Value = 2395190.4
$ws.Range("M5").Value = -1259
$ws.Range("N5").Value = -2395414.4
$ws.Range("H113").Value = 750.75
$ws.Range("I113").Value = 633
$ws.Range("J113").Value = 821.4
$ws.Range("K113").Value = 1899
$ws.Range("L113").Value = 2464.2
$ws.Range("M113").Value = 271
$ws.Range("N113").Value = -6804.2
$ws.Range("H117").Value = 1158.5555
$ws.Range("I117").Value = 485.8
$ws.Range("J117").Value = 1999.5
$ws.Range("K117").Value = 1457.4
$ws.Range("L117").Value = 5998.5
$ws.Range("M117").Value = 1984.6
$ws.Range("N117").Value = -12882.5
$ws.Range("H129").Value = 1988.8077
$ws.Range("I129").Value = 890
$ws.Range("J129").Value = 2570.5293
$ws.Range("K129").Value = 2670
$ws.Range("L129").Value = 7711.5879
$ws.Range("M129").Value = 2330
$ws.Range("N129").Value = -17711.5879
$ws.Range("H131").Value = 2494.6
$ws.Range("I131").Value = 561.46155
$ws.Range("J131").Value = 2820.974
$ws.Range("K131").Value = 1684.38465
$ws.Range("L131").Value = 8462.922
$ws.Range("M131").Value = 3355.61535
$ws.Range("N131").Value = -18542.922
$ws.Range("H135").Value = 439323.9
$ws.Range("I135").Value = 457
$ws.Range("J135").Value = 798396.8
$ws.Range("K135").Value = 4113
$ws.Range("L135").Value = 7185571.2
$ws.Range("M135").Value = -1578
$ws.Range("N135").Value = -7190641.2

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2991.111
$ws.Range("I102").Value = 1986.6666
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1986.6666
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -364.6666

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 883.3333
$ws.Range("I22").Value = 575
$ws.Range("J22").Value = 995.4545000000001
$ws.Range("K22").Value = 575
$ws.Range("L22").Value = 995.4545000000001
$ws.Range("M22").Value = -280
$ws.Range("N22").Value = -1585.4545
$ws.Range("H27").Value = 883.3333
$ws.Range("I27").Value = 575
$ws.Range("J27").Value = 995.4545000000001
$ws.Range("K27").Value = 575
$ws.Range("L27").Value = 995.4545000000001
$ws.Range("M27").Value = -468
$ws.Range("N27").Value = -1209.4545
$ws.Range("H55").Value = 353.89474
$ws.Range("I55").Value = 258.1111
$ws.Range("J55").Value = 440.1
$ws.Range("K55").Value = 258.1111
$ws.Range("L55").Value = 440.1
$ws.Range("M55").Value = -85.11110000000002
$ws.Range("N55").Value = -786.1
